# Updated symbol list on Thu Dec 29 03:32:35 UTC 2022 with GitHub Actions
#
# The sheet stores every data cell (Price/Volume/etc.) as literal TEXT,
# even though many of those strings look like numbers (e.g. "245.02").
# A plain `$cell.Value = "245.02"` would be auto-coerced by Excel into a
# real number (and can even pick up float rounding noise, e.g.
# 245.02000000000001), which would change the cell's type from string to
# number. To preserve the original text semantics we mark the cell as
# Text (NumberFormat "@") before assigning, then restore the default
# "Normal" style afterwards so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-PlainValue($range, [string]$value) {
    $ws.Range($range).Value = $value
}

# --- Price ticks (column D) for rows that only changed their price ---
Set-TextValue "D2"  "245.02"
Set-TextValue "D3"  "23.70"
Set-TextValue "D4"  "5.198"
Set-TextValue "D5"  "0.05748"
Set-TextValue "D6"  "6.460"
Set-TextValue "D7"  "3.248"
Set-TextValue "D8"  "0.8158"
Set-TextValue "D9"  "0.8676"

# --- Rows 10-18: coin list shifted up by one (wrapping "One" to the
#     bottom), each with its own refreshed price/volume-rank text ---
Set-PlainValue "B10" "WazirX"
Set-PlainValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue  "D10" "0.1378"
Set-PlainValue "E10" "9WazirXWRX"

Set-PlainValue "B11" "MandalaExchangeToken"
Set-PlainValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue  "D11" "0.06943"
Set-PlainValue "E11" "10MandalaExchangeTokenMDX"

Set-PlainValue "B12" "LiechtensteinCryptoassetsExchange"
Set-PlainValue "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue  "D12" "0.03232"
Set-PlainValue "E12" "11LiechtensteinCryptoassetsExchangeLCX"

Set-PlainValue "B13" "BitrueCoin"
Set-PlainValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue  "D13" "0.03017"
Set-PlainValue "E13" "12BitrueCoinBTR"

Set-PlainValue "B14" "BitMartToken"
Set-PlainValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue  "D14" "0.09329"
Set-PlainValue "E14" "13BitMartTokenBMX"

Set-PlainValue "B15" "MCDex"
Set-PlainValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue  "D15" "3.813"
Set-PlainValue "E15" "14MCDexMCB"

Set-PlainValue "B16" "BitForexToken"
Set-PlainValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue  "D16" "0.001520"
Set-PlainValue "E16" "15BitForexTokenBF"

Set-PlainValue "B17" "CoinExToken"
Set-PlainValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue  "D17" "0.04712"
Set-PlainValue "E17" "16CoinExTokenCET"

Set-PlainValue "B18" "One"
Set-PlainValue "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue  "D18" "0.0005976"
Set-PlainValue "E18" "17OneONE"

# --- Remaining standalone price ticks ---
Set-TextValue "D19" "0.006224"
Set-TextValue "D21" "0.004094"
Set-TextValue "D23" "3.575"
Set-TextValue "D24" "2.151"
Set-TextValue "D25" "0.3185"
Set-TextValue "D26" "0.1329"
Set-TextValue "D27" "0.0002327"
Set-TextValue "D40" "0.03711"
Set-TextValue "D42" "0.1050"

Set-TextValue  "D43" "0.002288"
Set-PlainValue "E43" "42CEJICEJIWorstin24h"

Set-TextValue "D44" "0.007019"
Set-TextValue "D45" "0.00005273"
Set-TextValue "D46" "0.00000000750"

Set-TextValue  "D47" "0.4298"
Set-PlainValue "E47" "46CoinbaseStockTokenCOIN"

Set-TextValue "D48" "0.002037"
Set-TextValue "D49" "0.00002099"
Set-TextValue "D50" "0.0001999"
